$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Characters(21,2).Text = "17"
$ws.Range("C9").Characters(27,9).Text = "4/21/2025"
$ws.Range("C9").Characters(47,9).Text = "4/27/2025"

# --- Cells whose value TYPE/STYLE changes (numeric <-> text placeholder) ---
# Copy from a same-shaped source cell so style index + literal type match exactly,
# then (for numeric targets) the copied value is already final; for text targets
# the copied shared-string value is already final too.
$ws.Range("C18").Copy($ws.Range("C15"))
$ws.Range("C14").Copy($ws.Range("C17"))
$ws.Range("C18").Copy($ws.Range("C27"))
$ws.Range("C18").Copy($ws.Range("C28"))
$ws.Range("C14").Copy($ws.Range("G29"))
$ws.Range("E14").Copy($ws.Range("H29"))
$ws.Range("C14").Copy($ws.Range("G30"))
$ws.Range("E14").Copy($ws.Range("H30"))
$ws.Range("C18").Copy($ws.Range("F33"))
$ws.Range("C18").Copy($ws.Range("I33"))

# --- Plain numeric value updates (style/type unchanged) ---
$ws.Range("F15").Value = 2
$ws.Range("I15").Value = 6
$ws.Range("K15").Value = 100
$ws.Range("L15").Value = 50
$ws.Range("M15").Value = 20
$ws.Range("N15").Value = -33.333333333333
$ws.Range("C16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 32
$ws.Range("J16").Value = 47
$ws.Range("K16").Value = -31.914893617021
$ws.Range("L16").Value = 10.344827586206
$ws.Range("M16").Value = -57.333333333333
$ws.Range("N16").Value = -84.236453201970
$ws.Range("E17").Value = -100
$ws.Range("F17").Value = 7
$ws.Range("H17").Value = -12.5
$ws.Range("I17").Value = 32
$ws.Range("J17").Value = 42
$ws.Range("K17").Value = -23.809523809523
$ws.Range("L17").Value = -39.622641509434
$ws.Range("M17").Value = -13.513513513513
$ws.Range("N17").Value = -63.636363636363
$ws.Range("F18").Value = 4
$ws.Range("H18").Value = 100
$ws.Range("I18").Value = 19
$ws.Range("K18").Value = 11.764705882352
$ws.Range("L18").Value = -40.625
$ws.Range("M18").Value = -78.160919540229
$ws.Range("N18").Value = -95.320197044335
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -16.666666666666
$ws.Range("F19").Value = 36
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = -25
$ws.Range("I19").Value = 159
$ws.Range("J19").Value = 191
$ws.Range("K19").Value = -16.753926701570
$ws.Range("L19").Value = -18.461538461538
$ws.Range("M19").Value = 3.246753246753
$ws.Range("N19").Value = -30.263157894736
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -20
$ws.Range("F20").Value = 12
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = -20
$ws.Range("I20").Value = 38
$ws.Range("J20").Value = 46
$ws.Range("K20").Value = -17.391304347826
$ws.Range("L20").Value = 11.764705882352
$ws.Range("M20").Value = -37.704918032786
$ws.Range("N20").Value = -95.927116827438
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = -14.285714285714
$ws.Range("G21").Value = 86
$ws.Range("H21").Value = -13.953488372093
$ws.Range("I21").Value = 286
$ws.Range("J21").Value = 346
$ws.Range("K21").Value = -17.341040462427
$ws.Range("L21").Value = -18.051575931232
$ws.Range("M21").Value = -31.742243436754
$ws.Range("N21").Value = -84.738527214514
$ws.Range("D23").Value = 2
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 13
$ws.Range("K23").Value = -69.230769230769
$ws.Range("C24").Value = 34
$ws.Range("D24").Value = 33
$ws.Range("E24").Value = 3.030303030303
$ws.Range("F24").Value = 108
$ws.Range("G24").Value = 114
$ws.Range("H24").Value = -5.263157894736
$ws.Range("I24").Value = 418
$ws.Range("J24").Value = 500
$ws.Range("K24").Value = -16.4
$ws.Range("L24").Value = 9.424083769633
$ws.Range("M24").Value = 44.137931034482
$ws.Range("C25").Value = 22
$ws.Range("D25").Value = 29
$ws.Range("E25").Value = -24.137931034482
$ws.Range("F25").Value = 82
$ws.Range("G25").Value = 100
$ws.Range("H25").Value = -18
$ws.Range("I25").Value = 284
$ws.Range("J25").Value = 417
$ws.Range("K25").Value = -31.89448441247
$ws.Range("L25").Value = 13.147410358565
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = -25
$ws.Range("F26").Value = 23
$ws.Range("G26").Value = 19
$ws.Range("H26").Value = 21.052631578947
$ws.Range("I26").Value = 86
$ws.Range("J26").Value = 92
$ws.Range("K26").Value = -6.521739130434
$ws.Range("L26").Value = 6.172839506172
$ws.Range("M26").Value = -9.473684210526
$ws.Range("F27").Value = 2
$ws.Range("I27").Value = 6
$ws.Range("K27").Value = 100
$ws.Range("L27").Value = -40
$ws.Range("E28").Value = 0
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = -66.666666666666
$ws.Range("I28").Value = 8
$ws.Range("J28").Value = 7
$ws.Range("K28").Value = 14.285714285714
$ws.Range("L28").Value = -33.333333333333
$ws.Range("N29").Value = -84.615384615384
$ws.Range("N30").Value = -83.333333333333
